$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Multi-Utilities(18)'
$ws.Range("B2").Value = 0.4701696426296317
$ws.Range("A3").Value = 'Mortgage Real Estate Investment Trust...(16)'
$ws.Range("B3").Value = 0.3724335807701659
$ws.Range("A4").Value = 'Road & Rail(22)'
$ws.Range("B4").Value = 0.3673881964468256
$ws.Range("A5").Value = 'Banks(251)'
$ws.Range("B5").Value = 0.3231659054417919
$ws.Range("A6").Value = 'Electric Utilities(28)'
$ws.Range("B6").Value = 0.3205586396304085
$ws.Range("A7").Value = 'Energy Equipment & Services(38)'
$ws.Range("B7").Value = 0.3194925571961612
$ws.Range("A8").Value = 'Marine(15)'
$ws.Range("B8").Value = 0.3070128247586821
$ws.Range("A9").Value = 'Auto Components(21)'
$ws.Range("B9").Value = 0.2938848054147093
$ws.Range("A10").Value = 'Machinery(86)'
$ws.Range("B10").Value = 0.2687891063247951
$ws.Range("A11").Value = 'Trading Companies & Distributors(25)'
$ws.Range("B11").Value = 0.2552831537526278
$ws.Range("A12").Value = 'Building Products(24)'
$ws.Range("B12").Value = 0.2491270427229211
$ws.Range("A13").Value = 'Specialty Retail(59)'
$ws.Range("B13").Value = 0.2488123605747006
$ws.Range("A14").Value = 'Construction & Engineering(21)'
$ws.Range("B14").Value = 0.2403076704806952
$ws.Range("A15").Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range("B15").Value = 0.2357986428000619
$ws.Range("A16").Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Range("B16").Value = 0.2332191071178296
$ws.Range("A17").Value = 'ETF(303)'
$ws.Range("B17").Value = 0.1968915353976214
$ws.Range("A18").Value = 'Hotels, Restaurants & Leisure(51)'
$ws.Range("B18").Value = 0.188945120984718
$ws.Range("A19").Value = 'Insurance(75)'
$ws.Range("B19").Value = 0.1840892825553289
$ws.Range("A20").Value = 'Capital Markets(76)'
$ws.Range("B20").Value = 0.1814315502615031
$ws.Range("A21").Value = 'Oil, Gas & Consumable Fuels(125)'
$ws.Range("B21").Value = 0.1775263443596889
$ws.Range("A22").Value = 'Semiconductors & Semiconductor Equipment(70)'
$ws.Range("B22").Value = 0.1734238930286855
$ws.Range("A23").Value = 'Professional Services(35)'
$ws.Range("B23").Value = 0.1672019193938927
$ws.Range("A24").Value = 'IT Services(52)'
$ws.Range("B24").Value = 0.1491692899065586
$ws.Range("A25").Value = 'Equity Real Estate Investment Trusts ...(98)'
$ws.Range("B25").Value = 0.1435237620395196
$ws.Range("A26").Value = 'Chemicals(52)'
$ws.Range("B26").Value = 0.1413109453507816
$ws.Range("A27").Value = 'Household Durables(39)'
$ws.Range("B27").Value = 0.1340439084648633
$ws.Range("A28").Value = 'Aerospace & Defense(37)'
$ws.Range("B28").Value = 0.1328854393022718
$ws.Range("A29").Value = 'Health Care Equipment & Supplies(86)'
$ws.Range("B29").Value = 0.1191024602879562
$ws.Range("A30").Value = 'Health Care Providers & Services(47)'
$ws.Range("B30").Value = 0.1143504074864952
$ws.Range("A31").Value = 'Commercial Services & Supplies(52)'
$ws.Range("B31").Value = 0.1076821198739539
$ws.Range("A32").Value = 'Metals & Mining(106)'
$ws.Range("B32").Value = 0.09393745275700958
$ws.Range("A33").Value = 'Electronic Equipment, Instruments & C...(78)'
$ws.Range("B33").Value = 0.07528335934453786
$ws.Range("A34").Value = 'Biotechnology(128)'
$ws.Range("B34").Value = 0.07421641189514432
$ws.Range("A35").Value = 'Software(70)'
$ws.Range("B35").Value = 0.07029839645929875
